$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B5: matching method description - update to mention Logit scale
$ws.Range("B5").Value = "1:1有放回匹配（Logit尺度）"

# B6: caliper value - originally stored as plain text, not a number.
# Force text formatting first so the numeric-looking string isn't
# auto-converted into a number by Excel's smart entry, then reset the
# cell style back to Normal so no extra "text number" style lingers.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "0.212436"
$ws.Range("B6").Style = "Normal"

# B9: matched sample count (numeric)
$ws.Range("B9").Value = 101

# B10: matching success rate - originally stored as plain text.
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "91.82%"
$ws.Range("B10").Style = "Normal"

# B11: post-matching treatment group sample count (numeric)
$ws.Range("B11").Value = 101

# B12: post-matching control group sample count (numeric)
$ws.Range("B12").Value = 101
